$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.940.59'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.671.22'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'214.92"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").Value = "'20.17"
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '1.906.79'
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").Value = '1.657.77'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = "'65.53"
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '26.940.27'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = "'233.97"
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("D19").Value = "'8.02"
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("D24").Value = "'2.18"
$ws.Range("D25").Value = "'145.89"
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("D26").Value = "'7.12"
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = "'15.96"
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("E28").Value = '  -1.30%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D33").Value = '1.461.44'
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("D36").Value = "'2.42"
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = "'0.580"
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("D39").Value = "'0.0170"
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").Value = "'1.06"
$ws.Range("E40").Value = '  +13.37%  '
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("D43").Value = "'2.29"
$ws.Range("E43").Value = '  +2.62%  '
$ws.Range("D44").Value = "'66.59"
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = '1.811.78'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("D46").Value = "'0.780"
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").Value = "'90.66"
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = '  +0.10%  '
